$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Grab the width of the neighbouring "In Advance" column (M) before inserting,
# so the freshly inserted column can be given the same width.
$neighbourWidth = $ws.Columns("M").ColumnWidth

# Make this the active sheet/tab (mirrors the workbook being saved with the
# "Repayment schedule" tab selected instead of "Transactions").
$ws.Activate()

# Inserting a whole column shifts the existing N/O/P columns (and their
# contents/widths) one place to the right, to O/P/Q.
$ws.Columns("N").Insert()

# New column keeps the same width as column M ("In Advance").
$ws.Columns("N").ColumnWidth = $neighbourWidth

# Restore/update the active selection on this sheet.
$ws.Range("J15").Select()
